$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (ID 14): new task description first, so the new shared string for the
# task text is registered before the "Påbörjad" status string below.
$ws.Range("C11").Value = "Jobba med Player-class (poängräkning, positionering med mera)"

# Row 5 (ID 8): status -> "Påbörjad", verklig tid -> 2
$ws.Range("D5").Value = "Påbörjad"
$ws.Range("F5").Value = 2

# Row 6 (ID 9): status -> "Påbörjad", skattad tid -> 3, verklig tid -> 1
$ws.Range("D6").Value = "Påbörjad"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1

# Row 11 (ID 14) continued: status -> "Påbörjad", skattad tid -> 4
$ws.Range("D11").Value = "Påbörjad"
$ws.Range("E11").Value = 4

# Leave the selection where the author left it before saving
$ws.Range("C9").Select()
